# "Add cantrals by cantons"
# The sheet had two header rows (E1/G1/I1/J1/K1 sparse labels + a second
# row with Hiver/Ete/Annee sub-labels). Replace that with a single,
# fully-populated header row (idx, idx2, Name, Date Start, Date End,
# (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year) and drop
# the old sub-label row, shifting the data rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old second header row (Hiver/Ete/Annee sub-labels) - data
# rows below it shift up to take its place.
$ws.Rows.Item(2).Delete()

# Rewrite row 1 as a single complete header row. Clear any leftover
# formatting from the old sparse header first so A1:E1 end up on the
# plain default style.
$ws.Range("A1:K1").ClearFormats()

$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 pick up the data font (Arial 9) like the rest of the table.
$hdr = $ws.Range("F1:K1")
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 9

# Match the selection left behind by the edit (row of new data, A2:K2).
[void]$ws.Range("A2:K2").Select()
